$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 461 (shifting rows 461:478 down to 462:479)
$ws.Rows("461:461").Insert(1)

# Populate the newly inserted row 461 with the new weekly price record.
$ws.Cells.Item(461, 1).Value = 4
$ws.Cells.Item(461, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(461, 3).Value = "Los Lagos"
$ws.Cells.Item(461, 4).Value = 45075
$ws.Cells.Item(461, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(461, 5).Value = 10
$ws.Cells.Item(461, 6).Value = 100112045
$ws.Cells.Item(461, 7).Value = "Zapallo"
$ws.Cells.Item(461, 8).Value = "Paine"
$ws.Cells.Item(461, 9).Value = "1a (guarda)"
$ws.Cells.Item(461, 10).Value = 250
$ws.Cells.Item(461, 11).Value = 550
$ws.Cells.Item(461, 12).Value = 550
$ws.Cells.Item(461, 13).Value = 550
$ws.Cells.Item(461, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(461, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(461, 16).Value = 550
$ws.Cells.Item(461, 17).Value = 1
$ws.Cells.Item(461, 18).Value = "Hortaliza"
